# Automatische test-sync: 2025-06-22 19:10:50
# Append a new row (40) of mail-log data to the "Logs" sheet and update
# the related dependent bits (conditional formatting ranges, Dashboard
# summary count) to match.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Logs")

# --- New row 40 content --------------------------------------------------
$ws.Range("A40").Value = "Technische storing"
$ws.Range("B40").Value = "mailmind.test@zohomail.eu"
$ws.Range("C40").Value = "De website werkt niet goed. Is hier iets mis mee?"
$ws.Range("D40").Value = "IT / Technisch probleem"
$ws.Range("E40").Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om dit verder te onderzoeken, heb ik wat meer informatie van u nodig. Kunt u specifiek aangeven welke problemen u ondervindt wanneer u de website bezoekt? Ziet u foutmeldingen, laadt de pagina niet, of zijn bepaalde functies niet beschikbaar? Alle details die u kunt geven zullen ons helpen het probleem op te lossen. `nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("F40").Value = "2025-06-22 19:10:15"
$ws.Range("G40").Value = "Ja"

# The multi-line text in E40 would otherwise trigger an automatic
# "customHeight" row resize; put the row back to the sheet's normal
# (non-custom) height, matching every other multi-line row already present.
$ws.Rows.Item(40).AutoFit()

# --- Extend the conditional formatting ranges to include the new row -----
$cfD = $ws.Range("D2:D39").FormatConditions.Item(1)
$cfD.ModifyAppliesToRange($ws.Range("D2:D40"))

$cfG = $ws.Range("G2:G39").FormatConditions.Item(1)
$cfG.ModifyAppliesToRange($ws.Range("G2:G40"))

# --- Update the Dashboard summary count for "IT / Technisch probleem" ----
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
